# Update Name of Algo
# Applies refreshed numeric results (RandomForest imputation) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 4.857600000000001
$ws.Range("A9").Value = -21.9738
$ws.Range("B12").Value = 5.641799999999998
$ws.Range("E13").Value = 16.75950000000001
$ws.Range("D15").Value = -8.285499999999999
$ws.Range("E16").Value = 16.09790000000001
$ws.Range("A18").Value = -22.20560000000001
$ws.Range("A20").Value = -19.13219999999999
$ws.Range("E20").Value = 16.03229999999999
$ws.Range("E24").Value = 16.61950000000001
$ws.Range("B26").Value = 4.034700000000004
$ws.Range("A27").Value = -21.47269999999996
$ws.Range("B27").Value = 4.882900000000004
$ws.Range("B29").Value = 5.098499999999999
$ws.Range("B37").Value = 9.105900000000004
$ws.Range("B38").Value = 4.484000000000002
$ws.Range("D38").Value = -8.937800000000005
$ws.Range("E39").Value = 16.0167
$ws.Range("D44").Value = -7.1649
$ws.Range("E48").Value = 17.5082
$ws.Range("B51").Value = 5.885699999999998
$ws.Range("D51").Value = -7.223799999999994
$ws.Range("E52").Value = 17.089
$ws.Range("B55").Value = 5.112399999999997
$ws.Range("E56").Value = 16.30280000000001
$ws.Range("D57").Value = -8.599600000000004
$ws.Range("D63").Value = -7.832299999999999
$ws.Range("A69").Value = -21.75959999999999
$ws.Range("B69").Value = 5.654399999999996
$ws.Range("B70").Value = 5.5921
$ws.Range("D70").Value = -6.913699999999999
$ws.Range("A76").Value = -19.27809999999998
$ws.Range("A82").Value = -22.02830000000001
$ws.Range("B83").Value = 5.649
$ws.Range("E84").Value = 16.77649999999999
$ws.Range("D99").Value = -8.018899999999999
$ws.Range("E100").Value = 16.45620000000001
$ws.Range("E101").Value = 16.79310000000001
$ws.Range("B102").Value = 8.647100000000007

$wb.Save()
